# Remove the "Ver no Jupiter..." and the "© 2020 ..." footer paragraphs,
# along with the blank paragraph that separates them from the preceding
# "LOB1004: Cálculo II (Requisito fraco)" line.
$d = $word.ActiveDocument

# Locate the paragraph that holds the "LOB1004..." text; the three
# paragraphs to delete are the ones immediately following it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOB1004: C*lculo II (Requisito fraco)*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $startPara = $target.Next()
    $endPara = $startPara.Next().Next()
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}
